$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2026-01-07 Wednesday", $false, $false, $false, $false, $false, $true, 1, $false, "2026-01-08 Thursday", 2) | Out-Null

# Update each equation cell in the table, in row-major order
$values = @(
    "0+75=",
    "75-53=",
    "73-48=",
    "76-14=",
    "1+4=",
    "47-22=",
    "20+18=",
    "22+44=",
    "64-52=",
    "96-4=",
    "21+0=",
    "64-48=",
    "25-17=",
    "43+12=",
    "8-7=",
    "74-43=",
    "68-24=",
    "23-0=",
    "10+19=",
    "33-8=",
    "24+4=",
    "26+32=",
    "22+20=",
    "91-73=",
    "52-15=",
    "91-66=",
    "40-35=",
    "91-36=",
    "64-23=",
    "59+19=",
    "27+45=",
    "49+11=",
    "90-51=",
    "9+68=",
    "65-27=",
    "71-4=",
    "1+73=",
    "14+13=",
    "62-41=",
    "92-25=",
    "92+2=",
    "82-71=",
    "26+72=",
    "41-27=",
    "86+5=",
    "19+57=",
    "40+31=",
    "41-16=",
    "20+53=",
    "70-23=",
    "96-0=",
    "42+40=",
    "86-23=",
    "95-64=",
    "14+11=",
    "18+73=",
    "59+37=",
    "93-89=",
    "86-83=",
    "37+6=",
    "74+23=",
    "16+52=",
    "21+52=",
    "3+30=",
    "92-53=",
    "5+91=",
    "78-11=",
    "1+22=",
    "63+20=",
    "25+16=",
    "80-75=",
    "60-33=",
    "28+27=",
    "65+34=",
    "85+5=",
    "26+59=",
    "73-67=",
    "33-23=",
    "63+25=",
    "57-22=",
    "19+43=",
    "66-31=",
    "2+21=",
    "57-38=",
    "18+26=",
    "58+22=",
    "59+39=",
    "89-50=",
    "63-39=",
    "88+2=",
    "49-42=",
    "50+24=",
    "60-45=",
    "59+27=",
    "79-43=",
    "98-69=",
    "66-22=",
    "81+15=",
    "32+66=",
    "33+17="
)

$t = $d.Tables.Item(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Host "Done. Updated $idx cells."